# ip_manager v3.7 import updates
$wb = $excel.ActiveWorkbook

# --- ip_address_list: mark "bewolktt" row as not favourite ---
$wsAddr = $wb.Worksheets.Item("ip_address_list")
$wsAddr.Range("E4").Value = 0

# --- ip_adress_fav_list: remove stale "bewolktt" favourite entry (row 3) ---
$wsFav = $wb.Worksheets.Item("ip_adress_fav_list")
$wsFav.Rows.Item(3).Delete()

# --- Settings: flip default startup window size flag ---
$wsSettings = $wb.Worksheets.Item("Settings")
$wsSettings.Range("B5").Value = 1
